$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 41714012
$ws.Range("I40").Value = 34333.332
$ws.Range("J40").Value = 55607236
$ws.Range("K40").Value = 34333.332
$ws.Range("L40").Value = 55607236
$ws.Range("M40").Value = -34158.332
$ws.Range("N40").Value = -55607586
$ws.Range("H86").Value = 2734136.5
$ws.Range("I86").Value = 3643707.2
$ws.Range("J86").Value = 5424.3335
$ws.Range("K86").Value = 3643707.2
$ws.Range("L86").Value = 5424.3335
$ws.Range("M86").Value = -3642584.2
$ws.Range("N86").Value = -7670.3335
$ws.Range("H89").Value = 2734136.5
$ws.Range("I89").Value = 3643707.2
$ws.Range("J89").Value = 5424.3335
$ws.Range("K89").Value = 18218536
$ws.Range("L89").Value = 27121.6675
$ws.Range("M89").Value = -18212920
$ws.Range("N89").Value = -38353.6675
$ws.Range("H131").Value = 1177.5
$ws.Range("I131").Value = 1212.2727
$ws.Range("J131").Value = 795
$ws.Range("K131").Value = 3636.8181
$ws.Range("L131").Value = 2385
$ws.Range("M131").Value = 1403.1819
$ws.Range("N131").Value = -12465
$ws.Range("H132").Value = 188567.08
$ws.Range("I132").Value = 820074.8
$ws.Range("K132").Value = 2460224.4
$ws.Range("M132").Value = -2457694.4
$ws.Range("H135").Value = 3233.875
$ws.Range("I135").Value = 838.9231
$ws.Range("J135").Value = 13612
$ws.Range("K135").Value = 7550.3079
$ws.Range("L135").Value = 122508
$ws.Range("M135").Value = -5015.3079
$ws.Range("N135").Value = -127578
$ws.Range("H138").Value = 3650.1638
$ws.Range("I138").Value = 1705
$ws.Range("J138").Value = 4284.4565
$ws.Range("K138").Value = 5115
$ws.Range("L138").Value = 12853.3695
$ws.Range("M138").Value = 25
$ws.Range("N138").Value = -23133.3695

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3597830
$ws.Range("I2").Value = 5823050.5
$ws.Range("K2").Value = 5823050.5
$ws.Range("M2").Value = -5822937.5
$ws.Range("H88").Value = 84966.164
$ws.Range("I88").Value = 1999
$ws.Range("J88").Value = 167933.33
$ws.Range("K88").Value = 1999
$ws.Range("L88").Value = 167933.33
$ws.Range("M88").Value = -1593
$ws.Range("N88").Value = -168745.33
$ws.Range("H91").Value = 84966.164
$ws.Range("I91").Value = 1999
$ws.Range("J91").Value = 167933.33
$ws.Range("K91").Value = 1999
$ws.Range("L91").Value = 167933.33
$ws.Range("M91").Value = -595
$ws.Range("N91").Value = -170741.33
$ws.Range("H116").Value = 3597830
$ws.Range("I116").Value = 5823050.5
$ws.Range("K116").Value = 5823050.5
$ws.Range("M116").Value = -5820756.5
$ws.Range("H122").Value = 6764.6816
$ws.Range("I122").Value = 3250.3333
$ws.Range("K122").Value = 9750.999899999999
$ws.Range("M122").Value = -7300.999899999999
$ws.Range("H131").Value = 70715
$ws.Range("J131").Value = 70715
$ws.Range("L131").Value = 70715
$ws.Range("N131").Value = -80795

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3597830
$ws.Range("I3").Value = 5823050.5
$ws.Range("K3").Value = 5823050.5
$ws.Range("M3").Value = -5822936.5
$ws.Range("H107").Value = 1387.6333
$ws.Range("I107").Value = 1660.5
$ws.Range("K107").Value = 1660.5
$ws.Range("M107").Value = 259.5
$ws.Range("H134").Value = 2143.7896
$ws.Range("I134").Value = 1782.8334
$ws.Range("K134").Value = 5348.5002
$ws.Range("M134").Value = -2813.5002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1774.8085
$ws.Range("I31").Value = 917.14703
$ws.Range("K31").Value = 917.14703
$ws.Range("M31").Value = -622.14703
$ws.Range("H34").Value = 1774.8085
$ws.Range("I34").Value = 917.14703
$ws.Range("K34").Value = 917.14703
$ws.Range("M34").Value = -715.14703
$ws.Range("H99").Value = 21001.25
$ws.Range("I99").Value = 56505.5
$ws.Range("K99").Value = 56505.5
$ws.Range("M99").Value = -55007.5
$ws.Range("H126").Value = 21001.25
$ws.Range("I126").Value = 56505.5
$ws.Range("K126").Value = 169516.5
$ws.Range("M126").Value = -167046.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H56").Value = 6899.091
$ws.Range("I56").Value = 6899.091
$ws.Range("K56").Value = 6899.091
$ws.Range("M56").Value = -6369.091
$ws.Range("H109").Value = 7323.75
$ws.Range("I109").Value = 1432
$ws.Range("J109").Value = 24999
$ws.Range("K109").Value = 4296
$ws.Range("L109").Value = 74997
$ws.Range("M109").Value = -3256
$ws.Range("N109").Value = -77077
$ws.Range("H113").Value = 989.6667
$ws.Range("I113").Value = 762.6667
$ws.Range("J113").Value = 1216.6666
$ws.Range("K113").Value = 2288.0001
$ws.Range("L113").Value = 3649.9998
$ws.Range("M113").Value = -118.0001000000002
$ws.Range("N113").Value = -7989.9998
$ws.Range("H114").Value = 3299.5
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 3299.5
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 9898.5
$ws.Range("M114").ClearContents()
$ws.Range("N114").Value = -16406.5
$ws.Range("H117").Value = 3687.6667
$ws.Range("J117").Value = 3687.6667
$ws.Range("L117").Value = 11063.0001
$ws.Range("N117").Value = -17947.0001
$ws.Range("H122").Value = 573
$ws.Range("J122").Value = 622.25
$ws.Range("L122").Value = 5600.25
$ws.Range("N122").Value = -10500.25
$ws.Range("H129").Value = 3548.75
$ws.Range("J129").Value = 4600
$ws.Range("L129").Value = 13800
$ws.Range("N129").Value = -23800
$ws.Range("H131").Value = 8476720
$ws.Range("J131").Value = 6946567.5
$ws.Range("L131").Value = 20839702.5
$ws.Range("N131").Value = -20849782.5
$ws.Range("H132").Value = 3398
$ws.Range("I132").Value = 2077.6
$ws.Range("K132").Value = 18698.4
$ws.Range("M132").Value = -16168.4
$ws.Range("H141").Value = 8557.842000000001
$ws.Range("I141").Value = 5739.615
$ws.Range("J141").Value = 14664
$ws.Range("K141").Value = 17218.845
$ws.Range("L141").Value = 43992
$ws.Range("M141").Value = -12038.845
$ws.Range("N141").Value = -54352

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 88955.53
$ws.Range("I80").Value = 203067.6
$ws.Range("K80").Value = 203067.6
$ws.Range("M80").Value = -202069.6
$ws.Range("H83").Value = 88955.53
$ws.Range("I83").Value = 203067.6
$ws.Range("K83").Value = 1015338
$ws.Range("M83").Value = -1010346
$ws.Range("H97").Value = 574.86664
$ws.Range("I97").Value = 631.63635
$ws.Range("K97").Value = 631.63635
$ws.Range("M97").Value = -135.63635
$ws.Range("H113").Value = 13249.75
$ws.Range("J113").Value = 16499.5
$ws.Range("L113").Value = 16499.5
$ws.Range("N113").Value = -20839.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6136.8335
$ws.Range("I7").Value = 4448.9
$ws.Range("K7").Value = 4448.9
$ws.Range("M7").Value = -4336.9
$ws.Range("H16").Value = 10002165
$ws.Range("I16").Value = 11112627
$ws.Range("J16").Value = 8001
$ws.Range("K16").Value = 11112627
$ws.Range("L16").Value = 8001
$ws.Range("M16").Value = -11112457
$ws.Range("N16").Value = -8341
$ws.Range("H22").Value = 1819.4667
$ws.Range("I22").Value = 779.6667
$ws.Range("K22").Value = 779.6667
$ws.Range("M22").Value = -484.6667
$ws.Range("H27").Value = 1819.4667
$ws.Range("I27").Value = 779.6667
$ws.Range("K27").Value = 779.6667
$ws.Range("M27").Value = -672.6667
$ws.Range("H100").Value = 3520.125
$ws.Range("I100").Value = 1624.25
$ws.Range("K100").Value = 1624.25
$ws.Range("M100").Value = -1083.25
$ws.Range("H126").Value = 6136.8335
$ws.Range("I126").Value = 4448.9
$ws.Range("K126").Value = 13346.7
$ws.Range("M126").Value = -10876.7

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H60").Value = 99999
$ws.Range("J60").Value = 99999
$ws.Range("L60").Value = 99999
$ws.Range("N60").Value = -101643
$ws.Range("H62").Value = 8317.333000000001
$ws.Range("I62").Value = 6634.6665
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 6634.6665
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -6010.6665
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 8317.333000000001
$ws.Range("I65").Value = 6634.6665
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 33173.3325
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -30053.3325
$ws.Range("N65").Value = -56240
$ws.Range("H107").Value = 2347.3333
$ws.Range("I107").Value = 2794.8
$ws.Range("J107").Value = 1068.8572
$ws.Range("K107").Value = 8384.400000000001
$ws.Range("L107").Value = 3206.5716
$ws.Range("M107").Value = -6464.400000000001
$ws.Range("N107").Value = -7046.571599999999
$ws.Range("H122").Value = 3928.38
$ws.Range("I122").Value = 3257.3784
$ws.Range("J122").Value = 5838.154
$ws.Range("K122").Value = 9772.135200000001
$ws.Range("L122").Value = 17514.462
$ws.Range("M122").Value = -7322.135200000001
$ws.Range("N122").Value = -22414.462
$ws.Range("H136").Value = 8053.222
$ws.Range("I136").Value = 3930.6843
$ws.Range("J136").Value = 9156.437
$ws.Range("K136").Value = 11792.0529
$ws.Range("L136").Value = 27469.311
$ws.Range("M136").Value = -9242.052899999999
$ws.Range("N136").Value = -32569.311
